# Updated cryptos list on Mon Jul 22 13:38:56 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# each coin row. Columns D/E are stored as plain text in the sheet, so
# values that look numeric (e.g. "0.999", "6.33") are entered with a
# leading apostrophe and the cell style is reset to "Normal" right after,
# which keeps them as literal text without leaving a numeric style behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.865.56'
$ws.Range("E2").Value = '  +1.68%  '

$ws.Range("D3").Value = '3.499.65'
$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = "'599.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.96%  '

$ws.Range("D6").Value = "'181.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.09%  '

$ws.Range("D7").Value = "'0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.18%  '

$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").Value = '3.500.09'
$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").Value = "'0.140"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.73%  '

$ws.Range("D11").Value = "'7.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.16%  '

$ws.Range("D12").Value = "'0.437"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.76%  '

$ws.Range("D13").Value = '4.094.51'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").Value = "'32.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.64%  '

$ws.Range("D15").Value = "'0.135"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.40%  '

$ws.Range("D16").Value = '67.611.92'
$ws.Range("E16").Value = '  +1.20%  '

$ws.Range("D17").Value = "'0.0000179"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.83%  '

$ws.Range("D18").Value = '3.492.80'
$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("D19").Value = "'6.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.33%  '

$ws.Range("D20").Value = "'14.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.29%  '

$ws.Range("D21").Value = "'392.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.69%  '

$ws.Range("D22").Value = "'8.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.09%  '

$ws.Range("D23").Value = "'73.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.21%  '

$ws.Range("D24").Value = "'0.543"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.81%  '

$ws.Range("E25").Value = '  +0.55%  '

$ws.Range("D26").Value = "'5.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.35%  '

$ws.Range("D27").Value = "'0.0000124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.70%  '

$ws.Range("D28").Value = "'10.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.70%  '

$ws.Range("E29").Value = '  -1.95%  '

$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").Value = "'6.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.18%  '

$ws.Range("D32").Value = "'1.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.09%  '

$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("D34").Value = "'23.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("E35").Value = '  +0.93%  '

$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").Value = "'1.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.36%  '

$ws.Range("D38").Value = "'162.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.85%  '

$ws.Range("D39").Value = "'0.888"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.66%  '

$ws.Range("D40").Value = "'2.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.59%  '

$ws.Range("D41").Value = "'1.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("E42").Value = '  +0.07%  '

$ws.Range("D43").Value = "'4.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.48%  '

$ws.Range("D44").Value = '2.833.46'
$ws.Range("E44").Value = '  +0.11%  '

$ws.Range("D45").Value = "'26.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.77%  '

$ws.Range("D46").Value = "'26.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.63%  '

$ws.Range("D47").Value = "'0.0729"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.25%  '

$ws.Range("D48").Value = "'41.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.07%  '

$ws.Range("D49").Value = "'0.0302"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("D50").Value = "'335.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.21%  '

$ws.Range("E51").Value = '  -0.93%  '
